$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column C (course_x), shifting existing
# course_x..instructor_y columns from C:J to E:L.
$ws.Range("C:D").Insert()

# New header cells for the inserted columns, matching the bold style
# used by the rest of the header row.
$ws.Range("C1").Value = "image_x"
$ws.Range("D1").Value = "image_y"
$ws.Range("C1:D1").Font.Bold = $true

$ws.Range("C2").Value = 1440
$ws.Range("D2").Value = 990

$ws.Range("C3").Value = 1440
$ws.Range("D3").Value = 990

# Size the new columns similarly to the other bestFit columns on this
# sheet (closest width the engine's column-width quantization can reach).
$ws.Columns("C").ColumnWidth = 6.666666666666666
$ws.Columns("D").ColumnWidth = 6.666666666666666

# Move the active selection (no longer the header row) as in the saved file.
$ws.Range("D6").Select()
